# edit.ps1 - apply the workbook update described by the diff:
#  - fill in rows 6-11 of "Таблица2" (previously blank placeholder rows) with
#    new trip records
#  - fix a data-entry mistake on row 4 (price tier id/from/to)
#  - fix the "id цен диапазона" sequence (G18:G20) in the lookup table
#  - move the active selection to E22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Copy the cell formatting of existing filled rows down onto the blank
#    rows so the new rows look like their siblings (fills/borders/number
#    formats for currency columns etc.).
# ---------------------------------------------------------------------------
$ws.Range("A3:K3").Copy() | Out-Null
$ws.Range("A6:K6").PasteSpecial(-4122) | Out-Null

$ws.Range("A2:K2").Copy() | Out-Null
$ws.Range("A7:K7").PasteSpecial(-4122) | Out-Null
$ws.Range("A8:K8").PasteSpecial(-4122) | Out-Null

$ws.Range("A4:K4").Copy() | Out-Null
$ws.Range("A9:K9").PasteSpecial(-4122) | Out-Null

$ws.Range("A2:K2").Copy() | Out-Null
$ws.Range("A10:K10").PasteSpecial(-4122) | Out-Null
$ws.Range("A11:K11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Fill in the new trip rows 6-11.
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Москва"
$ws.Range("B6").Value = "Минск"
$ws.Range("C6").Value = "22.04.2023"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "Самолет"
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 5000
$ws.Range("J6").Value = 10000

$ws.Range("A7").Value = "Сургут "
$ws.Range("B7").Value = "Санкт-Петербург"
$ws.Range("C7").Value = "23.04.2023"
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "Поезд"
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 5000

$ws.Range("A8").Value = "Санкт-Петербург"
$ws.Range("B8").Value = "Минск"
$ws.Range("C8").Value = "23.04.2023"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = "Самолет"
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 5000

$ws.Range("A9").Value = "Сургут"
$ws.Range("B9").Value = "Минск"
$ws.Range("C9").Value = "23.04.2023"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "Самолет"
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = "Санкт-Петербург"
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = 5000
$ws.Range("J9").Value = 10000

$ws.Range("A10").Value = "Сургут "
$ws.Range("B10").Value = "Санкт-Петербург"
$ws.Range("C10").Value = "24.04.2023"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = "Самолет"
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1000
$ws.Range("J10").Value = 5000

$ws.Range("A11").Value = "Санкт-Петербург"
$ws.Range("B11").Value = "Минск"
$ws.Range("C11").Value = "23.04.2023"
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = "Поезд"
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1000
$ws.Range("J11").Value = 5000

# ---------------------------------------------------------------------------
# 3) Correct row 4's price tier (id цен диапазона 1 -> 2, от/до 1000/5000 ->
#    5000/10000) to match the "Поезд" tier actually used.
# ---------------------------------------------------------------------------
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 10000

# ---------------------------------------------------------------------------
# 4) Fix the "id цен диапазона" reference column (G18:G20) which was
#    mistakenly left at 1 for every row instead of being sequential.
# ---------------------------------------------------------------------------
$ws.Range("G18").Value = 2
$ws.Range("G19").Value = 3
$ws.Range("G20").Value = 4

# ---------------------------------------------------------------------------
# 5) Move the active selection, matching the saved cursor position.
# ---------------------------------------------------------------------------
$ws.Range("E22").Select() | Out-Null
